$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff" — refresh the handoff/generate timestamps
# for the f05f1d23-b45b-448a-a6af-01d0dafae906.md row (row 7 in every table)
# to reflect a freshly generated handoff package.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-27 08:41:17"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-27 08:41:12"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-27 08:41:17"
